# Commit: "load ue5 nav mesh bin"
# The underlying change renames the worksheet from "scene" to "mainscene"
# (matching the workbook file name mainscene.xlsx) and updates the
# active-cell selection that had been left on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet; this also updates the <sheet name="..."> entry.
$ws.Name = "mainscene"

# The workbook's hidden _FilterDatabase defined name is qualified with the
# sheet name (it refers to an invalid #REF! range). Renaming the sheet
# drops the qualifier because the target range no longer resolves, so
# restore the sheet-qualified text explicitly.
foreach ($dn in $wb.Names) {
    if ($dn.Name -like "*_FilterDatabase*") {
        $dn.RefersTo = "=mainscene!#REF!"
    }
}

# Move the selection/active cell from B8:B23 to the single cell E8.
$null = $ws.Range("E8").Select()
